$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.208.44"
Set-TextValue "E2" "  -4.93%  "
Set-TextValue "D3" "3.253.39"
Set-TextValue "E3" "  -7.58%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "597.16"
Set-TextValue "E5" "  -3.53%  "
Set-TextValue "D6" "150.45"
Set-TextValue "E6" "  -13.06%  "
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "3.244.16"
Set-TextValue "E8" "  -7.73%  "
Set-TextValue "D9" "0.541"
Set-TextValue "E9" "  -11.42%  "
Set-TextValue "D10" "0.171"
Set-TextValue "E10" "  -13.64%  "
Set-TextValue "E11" "  -5.31%  "
Set-TextValue "D12" "0.505"
Set-TextValue "E12" "  -13.98%  "
Set-TextValue "D13" "38.01"
Set-TextValue "E13" "  -18.11%  "
Set-TextValue "E14" "  -12.29%  "
Set-TextValue "D15" "3.774.26"
Set-TextValue "E15" "  -7.74%  "
Set-TextValue "D16" "67.252.92"
Set-TextValue "E16" "  -4.99%  "
Set-TextValue "D17" "3.260.48"
Set-TextValue "E17" "  -7.54%  "
Set-TextValue "D18" "540.53"
Set-TextValue "E18" "  -11.20%  "
Set-TextValue "E19" "  -6.16%  "
Set-TextValue "D20" "7.19"
Set-TextValue "E20" "  -14.19%  "
Set-TextValue "D21" "15.09"
Set-TextValue "E21" "  -14.90%  "
Set-TextValue "D22" "0.759"
Set-TextValue "E22" "  -14.07%  "
Set-TextValue "D23" "7.86"
Set-TextValue "E23" "  -14.55%  "
Set-TextValue "D24" "85.37"
Set-TextValue "E24" "  -12.74%  "
Set-TextValue "D25" "13.43"
Set-TextValue "E25" "  -13.98%  "
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  -0.08%  "
Set-TextValue "D27" "3.26"
Set-TextValue "E27" "  -12.36%  "
Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "29.29"
Set-TextValue "E28" "  -12.98%  "
Set-TextValue "B29" "RenderToken"
Set-TextValue "C29" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D29" "8.01"
Set-TextValue "E29" "  -11.63%  "
Set-TextValue "D30" "2.12"
Set-TextValue "E30" "  -17.47%  "
Set-TextValue "D31" "2.66"
Set-TextValue "E31" "  -11.24%  "
Set-TextValue "E32" "  -12.43%  "
Set-TextValue "B33" "Bittensor"
Set-TextValue "C33" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D33" "542.31"
Set-TextValue "E33" "  -15.42%  "
Set-TextValue "B34" "Filecoin"
Set-TextValue "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "6.63"
Set-TextValue "E34" "  -18.07%  "
Set-TextValue "D35" "5.69"
Set-TextValue "E35" "  -16.63%  "
Set-TextValue "E36" "  +0.12%  "
Set-TextValue "D37" "0.0441"
Set-TextValue "E37" "  -9.27%  "
Set-TextValue "D38" "53.12"
Set-TextValue "E38" "  -6.36%  "
Set-TextValue "D39" "0.0851"
Set-TextValue "E39" "  -14.65%  "
Set-TextValue "D40" "9.14"
Set-TextValue "E40" "  -15.42%  "
Set-TextValue "D41" "0.128"
Set-TextValue "E41" "  -10.03%  "
Set-TextValue "D42" "2.920.14"
Set-TextValue "E42" "  -12.88%  "
Set-TextValue "E43" "  -22.53%  "
Set-TextValue "D44" "0.261"
Set-TextValue "E44" "  -16.29%  "
Set-TextValue "D45" "0.0₃0582"
Set-TextValue "E45" "  -19.03%  "
Set-TextValue "D46" "2.17"
Set-TextValue "E46" "  -14.47%  "
Set-TextValue "D47" "26.45"
Set-TextValue "E47" "  -16.84%  "
Set-TextValue "E48" "  -0.07%  "
Set-TextValue "D49" "127.43"
Set-TextValue "E49" "  -5.27%  "
Set-TextValue "D50" "2.33"
Set-TextValue "E50" "  -21.20%  "
Set-TextValue "D51" "0.113"
Set-TextValue "E51" "  -12.79%  "
